$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55 (shifts existing rows 55:101 down to 56:102)
$ws.Rows.Item(55).Insert()

# New record inserted at row 55
$ws.Range("A55").Value = 11
$ws.Range("B55").Value = "Vega Monumental Concepción"
$ws.Range("C55").Value = "Bíobío"
$ws.Range("D55").Value = 44586
$ws.Range("D55").NumberFormat = $ws.Range("D56").NumberFormat
$ws.Range("E55").Value = 8
$ws.Range("F55").Value = 100112043
$ws.Range("G55").Value = "Pepino ensalada"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 90
$ws.Range("K55").Value = 11000
$ws.Range("L55").Value = 12000
$ws.Range("M55").Value = 11444
$ws.Range("N55").Value = '$/caja 60 unidades'
$ws.Range("O55").Value = "Región de Arica y Parinacota"
$ws.Range("P55").Value = 191
$ws.Range("Q55").Value = 60
$ws.Range("R55").Value = "Hortaliza"
